# ---------------------------------------------------------------------------
# Adds results from the "new optimization" run: renames the two worksheets,
# wires up a second (dataOutput) named range / text import next to the
# existing detailedOutput one, and populates "neue Optimierung" (formerly
# Tabelle2) with the imported table (rows 1-5) plus the "Imported fields"
# summary block (rows 17-25). Also fixes up the chart series that pointed at
# the renamed "Tabelle1" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOld = $wb.Worksheets.Item(1)   # Tabelle1       -> altes Zeug
$wsNew = $wb.Worksheets.Item(2)   # Tabelle2       -> neue Optimierung

# ---------------------------------------------------------------------------
# 1. Rename the sheets
# ---------------------------------------------------------------------------
$wsOld.Name = "altes Zeug"
$wsNew.Name = "neue Optimierung"

# ---------------------------------------------------------------------------
# 2. Fix up the three charts on "altes Zeug" that referenced the old sheet
#    name explicitly inside their SERIES() formula (chart1 -> col F,
#    chart2 -> col D, chart3 -> col G).
# ---------------------------------------------------------------------------
$co1 = $wsOld.ChartObjects().Item(1)
$co1.Chart.SeriesCollection().Item(1).Formula = "=SERIES(,'altes Zeug'!`$A`$2:`$A`$9,'altes Zeug'!`$F`$2:`$F`$9,1)"

$co2 = $wsOld.ChartObjects().Item(2)
$co2.Chart.SeriesCollection().Item(1).Formula = "=SERIES(,'altes Zeug'!`$A`$2:`$A`$9,'altes Zeug'!`$D`$2:`$D`$9,1)"

$co3 = $wsOld.ChartObjects().Item(3)
$co3.Chart.SeriesCollection().Item(1).Formula = "=SERIES(,'altes Zeug'!`$A`$2:`$A`$9,'altes Zeug'!`$G`$2:`$G`$9,1)"

# ---------------------------------------------------------------------------
# 3. New defined name for the imported text data on "neue Optimierung"
#    (mirrors the existing detailedOutput name on "altes Zeug").
# ---------------------------------------------------------------------------
$wsNew.Names.Add("dataOutput", "='neue Optimierung'!`$A`$18:`$B`$25")

# ---------------------------------------------------------------------------
# 4. Column widths / layout for "neue Optimierung" (matches the imported
#    table layout already used on "altes Zeug").
# ---------------------------------------------------------------------------
$colWidths = @(18.5703125, 12, 11.28515625, 13.85546875, 12.85546875, 22.7109375, 12.140625, 15.5703125, 14.85546875, 16.28515625, 17.140625, 20.28515625, 14.85546875, 16.28515625, 17.42578125, 19.5703125, 20.7109375, 24.140625)
for ($c = 1; $c -le $colWidths.Length; $c++) {
    $wsNew.Columns.Item($c).ColumnWidth = $colWidths[$c - 1] - 0.8333333333333334
}

# ---------------------------------------------------------------------------
# 5. Header row (row 1)
# ---------------------------------------------------------------------------
$headers = @("City", "Energy cost in €/kWh", "Penalty in €", "Vstorage in m³", "Acollector m²", "thicknessInsulation in m", "solarfraction", "costfunction in €", "costStorage in €", "costCollector in €", "costInsulation in €", "costHeaterEnergy in €", "costPenalty in €", "heaterEnergy in J", "radiatorEnergy in J", "heaterEnergy in kWh", "radiatorEnergy in kWh", "execution time in GenOpt")
for ($c = 1; $c -le $headers.Length; $c++) {
    $wsNew.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ---------------------------------------------------------------------------
# 6. Imported data table (rows 2-5)
# ---------------------------------------------------------------------------
$wsNew.Range("A2").Value = "Chicago (pen = 0 €, EP = 0.08 €)"
$wsNew.Range("B2").Value = 0.08
$wsNew.Range("C2").Value = 0
$wsNew.Range("D2").Value = 1
$wsNew.Range("E2").Value = 3.5
$wsNew.Range("F2").Value = 0.15375
$wsNew.Range("G2").Value = 0.0845623325557276
$wsNew.Range("H2").Value = 1572.73
$wsNew.Range("I2").Value = 82.4904999999999
$wsNew.Range("J2").Value = 35
$wsNew.Range("K2").Value = 907.344621333333
$wsNew.Range("L2").Value = 547.89584913244
$wsNew.Range("M2").Value = 0
$wsNew.Range("N2").Value = 24655313210.9598
$wsNew.Range("O2").Value = 26932814857.6109
$wsNew.Range("P2").Formula = "=N2/3600/1000"
$wsNew.Range("Q2").Formula = "=O2/3600/1000"
$wsNew.Range("R2").Value = 0.031655092592592596
$wsNew.Range("R2").NumberFormat = "h:mm:ss"

$wsNew.Range("A3").Value = "Chicago (pen = 1500 €, EP = 0.08 €)"
$wsNew.Range("B3").Value = 0.08
$wsNew.Range("C3").Value = 1500
$wsNew.Range("D3").Value = 29.5
$wsNew.Range("E3").Value = 9.5
$wsNew.Range("F3").Value = 0.19875
$wsNew.Range("G3").Value = 0.500448289933642
$wsNew.Range("H3").Value = 1856.28
$wsNew.Range("I3").Value = 506.090909578394
$wsNew.Range("J3").Value = 95
$wsNew.Range("K3").Value = 986.925837333333
$wsNew.Range("L3").Value = 268.262182913226
$wsNew.Range("M3").Value = 0
$wsNew.Range("N3").Value = 12071798231.0952
$wsNew.Range("O3").Value = 24165262550.0804
$wsNew.Range("P3").Formula = "=N3/3600/1000"
$wsNew.Range("Q3").Formula = "=O3/3600/1000"
$wsNew.Range("R3").Value = 0.03542824074074074
$wsNew.Range("R3").NumberFormat = "h:mm:ss"

$wsNew.Range("A4").Value = "San Francisco (pen = 0 €, EP = 0.08 €)"
$wsNew.Range("B4").Value = 0.08
$wsNew.Range("C4").Value = 0
$wsNew.Range("D4").Value = 1
$wsNew.Range("E4").Value = 1
$wsNew.Range("F4").Value = 0.06375
$wsNew.Range("G4").Value = 0.196343803338515
$wsNew.Range("H4").Value = 905.79
$wsNew.Range("I4").Value = 82.4904999999999
$wsNew.Range("J4").Value = 10
$wsNew.Range("K4").Value = 748.182189333333
$wsNew.Range("L4").Value = 65.1204593271224
$wsNew.Range("M4").Value = 0
$wsNew.Range("N4").Value = 2930420669.7205
$wsNew.Range("O4").Value = 3646361070.68413
$wsNew.Range("P4").Formula = "=N4/3600/1000"
$wsNew.Range("Q4").Formula = "=O4/3600/1000"
$wsNew.Range("R4").Value = 0.018217592592592594
$wsNew.Range("R4").NumberFormat = "h:mm:ss"

$wsNew.Range("A5").Value = "San Francisco (pen = 1500 €, EP = 0.08 €)"
$wsNew.Range("B5").Value = 0.08
$wsNew.Range("C5").Value = 1500
$wsNew.Range("D5").Value = 1
$wsNew.Range("E5").Value = 6.5
$wsNew.Range("F5").Value = 0.06
$wsNew.Range("G5").Value = 0.506218821787891
$wsNew.Range("H5").Value = 932.78
$wsNew.Range("I5").Value = 82.4904999999999
$wsNew.Range("J5").Value = 65
$wsNew.Range("K5").Value = 741.550421333333
$wsNew.Range("L5").Value = 43.7426394853869
$wsNew.Range("M5").Value = 0
$wsNew.Range("N5").Value = 1968418776.84241
$wsNew.Range("O5").Value = 3986419215.02495
$wsNew.Range("P5").Formula = "=N5/3600/1000"
$wsNew.Range("Q5").Formula = "=O5/3600/1000"
$wsNew.Range("R5").Value = 0.01765046296296296
$wsNew.Range("R5").NumberFormat = "h:mm:ss"

# ---------------------------------------------------------------------------
# 7. "Imported fields" summary block (rows 18-25), header at row 17 goes in
#    last so that the shared-string table order matches the append order
#    used by the model (rows 2,4 then 18-25 then 17).
# ---------------------------------------------------------------------------
$wsNew.Range("A18").Value = "costStorage = "
$wsNew.Range("B18").Value = 82.4904999999999

$wsNew.Range("A19").Value = "costCollector = "
$wsNew.Range("B19").Value = 65

$wsNew.Range("A20").Value = "costInsulation = "
$wsNew.Range("B20").Value = 741.550421333333

$wsNew.Range("A21").Value = "heaterEnergy = "
$wsNew.Range("B21").Value = 1968418776.84241

$wsNew.Range("A22").Value = "radiatorEnergy = "
$wsNew.Range("B22").Value = 3986419215.02495

$wsNew.Range("A23").Value = "solarfraction = "
$wsNew.Range("B23").Value = 0.506218821787891

$wsNew.Range("A24").Value = "costHeaterEnergy = "
$wsNew.Range("B24").Value = 43.7426394853869

$wsNew.Range("A25").Value = "costPenalty = "
$wsNew.Range("B25").Value = 0

$wsNew.Range("A17").Value = "Imported fields"

# ---------------------------------------------------------------------------
# 8. Sheet view / selection state (matches the post-edit workbook: the new
#    sheet becomes the active tab with A18 selected, the old sheet loses its
#    previous scroll position / selection).
# ---------------------------------------------------------------------------
$wsOld.Select()
$wsOld.Range("A2:A9").Select()
$wsNew.Select()
$wsNew.Range("A18").Select()
